# Apply 1.5 line spacing (w:line="360" w:lineRule="auto") to every
# paragraph in the document body, preserving any spacing already set
# (e.g. w:after="53"/"55"). This mirrors Word's "Line Spacing: 1.5
# lines" applied to Ctrl+A-selected content.

$d = $word.ActiveDocument

$d.Content.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceMultiple
$d.Content.ParagraphFormat.LineSpacing = 18      # 18pt => 360 twips => 1.5 lines
